# "Added cards to every user" — fill in the missing Card # (column E)
# values for every user row on the UserList sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserList")
$ws.Activate()

$cardNumbers = @{
    3 = 100000001
    4 = 100000002
    5 = 100000003
    6 = 100000004
    7 = 100000005
    8 = 100000006
    9 = 100000007
}

foreach ($row in $cardNumbers.Keys) {
    $ws.Cells.Item($row, 5).Value = $cardNumbers[$row]
}

$ws.Range("E9").Select()
